# issue #5: stock data output to json file
#
# The "股票" (Stock) worksheet gains a new "property_category" column
# (value "stock" for every data row), inserted immediately before the
# existing "date" column. Everything that used to live at/after that
# column (date, legislator_name, legislator_id) shifts one column to
# the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at H (before the current "date" column), pushing
# date / legislator_name / legislator_id one column to the right.
$ws.Columns.Item(8).Insert()

# New header cell for the inserted column.
$ws.Range("H1").Value = "property_category"

# New data values: every stock row is categorized as "stock".
$ws.Range("H2:H7").Value = "stock"
